$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "PR for Config_Child_Windows_20220521." becomes two
#    runs: "There is a "ding" when creating a new folder and a new folder
#    group." followed by " Enter button."
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "There is a $([char]0x201C)ding$([char]0x201D) when creating a new folder and a new folder group."

# Type the second sentence as a separate paragraph first, then remove the
# paragraph mark that divides them. This keeps the two sentences as two
# distinct <w:r> runs instead of them being coalesced into a single run.
$p1.Range.InsertParagraphAfter()
$p1Tail = $d.Paragraphs(2)
$p1Tail.Range.Text = " Enter button."
$d.Range($p1.Range.End - 1, $p1.Range.End).Delete()

# ---------------------------------------------------------------------
# 2) Second paragraph ("The textbox is to short. Not very practical when
#    first rendered.") is cleared out, leaving an empty paragraph behind.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$d.Range($p2.Range.Start, $p2.Range.End - 1).Delete()

# ---------------------------------------------------------------------
# 3) The "Detect backspace and arrows?" paragraph is removed completely,
#    including the paragraph mark that precedes it.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Detect backspace and arrows?*") {
        $prev = $p.Previous()
        $d.Range($prev.Range.End, $p.Range.End).Delete()
        break
    }
}

foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [$($p.Range.Text)]"
}
